# Apply "Ran code for averaged intensities on spiral schemes" update.
#
# Semantics of the change:
#  - Three new sampling schemes were run and their averaged-intensity rows were
#    inserted into the table right after "Ring Perpendicular to TD" and before
#    "NoRotation-tilt60deg":
#        Gaussian-Quadrature              (already existed - now recomputed)
#        Spiral-90deg-10rot-5space        (new)
#        Spiral-90deg-15rot-5space        (new)
#        Spiral-90deg-10rot-3space        (new)
#  - The rows that used to follow ("NoRotation-tilt60deg" ... "HexGrid-60degTilt5degRes")
#    shift down by 4 data rows (worksheet rows 10-15 -> 14-19), and the
#    "Gaussian-Quadrature" row (previously the very last row, worksheet row 16)
#    moves up to become worksheet row 10 (directly after "Ring Perpendicular to TD").
#  - The net result is 3 brand-new rows (11, 12, 13) and the table grows from
#    A1:P16 to A1:P19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Gaussian-Quadrature --------------------------------------------
$ws.Range("B10").Value2 = "Gaussian-Quadrature"
$ws.Range("C10").Value2 = 2.018785056086185
$ws.Range("D10").Value2 = 2.056121120826005
$ws.Range("E10").Value2 = 1.634508309208395
$ws.Range("F10").Value2 = 0.7388292315594466
$ws.Range("G10").Value2 = 2.018785056086185
$ws.Range("H10").Value2 = 2.056121120826005
$ws.Range("I10").Value2 = 1.050463488344989
$ws.Range("J10").Value2 = 0.5806850065064968
$ws.Range("K10").Value2 = 1.021648733139336
$ws.Range("L10").Value2 = 0.8781271138853387
$ws.Range("M10").Value2 = 2.018785056086185
$ws.Range("N10").Value2 = 1.8453147150172
$ws.Range("O10").Value2 = 1.612060929420008
$ws.Range("P10").Value2 = 1.247396007444524

# --- Row 11: Spiral-90deg-10rot-5space (new) --------------------------------
$ws.Range("B11").Value2 = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value2 = 0
$ws.Range("D11").Value2 = 3.563089281232255
$ws.Range("E11").Value2 = 1.559426470653133
$ws.Range("F11").Value2 = 0.5167963296042508
$ws.Range("G11").Value2 = 0
$ws.Range("H11").Value2 = 3.563089281232255
$ws.Range("I11").Value2 = 0.4979775931499223
$ws.Range("J11").Value2 = 1.78345674844514
$ws.Range("K11").Value2 = 0.08727927503991484
$ws.Range("L11").Value2 = 1.661787810905623
$ws.Range("M11").Value2 = 0
$ws.Range("N11").Value2 = 2.561257875942694
$ws.Range("O11").Value2 = 1.40982802037241
$ws.Range("P11").Value2 = 1.20872668862878

# --- Row 12: Spiral-90deg-15rot-5space (new) --------------------------------
$ws.Range("B12").Value2 = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value2 = 0
$ws.Range("D12").Value2 = 3.57031561225003
$ws.Range("E12").Value2 = 1.540973837564633
$ws.Range("F12").Value2 = 0.5182766746288745
$ws.Range("G12").Value2 = 0
$ws.Range("H12").Value2 = 3.57031561225003
$ws.Range("I12").Value2 = 0.4913381217411473
$ws.Range("J12").Value2 = 1.787669925213792
$ws.Range("K12").Value2 = 0.0874534275581257
$ws.Range("L12").Value2 = 1.665985788892413
$ws.Range("M12").Value2 = 0
$ws.Range("N12").Value2 = 2.555644724907332
$ws.Range("O12").Value2 = 1.407391531110884
$ws.Range("P12").Value2 = 1.207751673481127

# --- Row 13: Spiral-90deg-10rot-3space (new) --------------------------------
$ws.Range("B13").Value2 = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value2 = 0
$ws.Range("D13").Value2 = 3.562501612822962
$ws.Range("E13").Value2 = 1.553384977716487
$ws.Range("F13").Value2 = 0.5172521053421978
$ws.Range("G13").Value2 = 0
$ws.Range("H13").Value2 = 3.562501612822962
$ws.Range("I13").Value2 = 0.4953156636663871
$ws.Range("J13").Value2 = 1.786140249288871
$ws.Range("K13").Value2 = 0.08738987830544981
$ws.Range("L13").Value2 = 1.662430743069581
$ws.Range("M13").Value2 = 0
$ws.Range("N13").Value2 = 2.557943295269724
$ws.Range("O13").Value2 = 1.408284673970412
$ws.Range("P13").Value2 = 1.208051903776492

# --- Row 14: NoRotation-tilt60deg (shifted from old row 10) ----------------
$ws.Range("B14").Value2 = "NoRotation-tilt60deg"
$ws.Range("C14").Value2 = 0
$ws.Range("D14").Value2 = 2.392376000000007
$ws.Range("E14").Value2 = 4.576016000000003
$ws.Range("F14").Value2 = 0.2558239999999998
$ws.Range("G14").Value2 = 0
$ws.Range("H14").Value2 = 2.392376000000007
$ws.Range("I14").Value2 = 1.39418
$ws.Range("J14").Value2 = 1.170812
$ws.Range("K14").Value2 = 0.05374399999999981
$ws.Range("L14").Value2 = 0.979444000000003
$ws.Range("M14").Value2 = 0
$ws.Range("N14").Value2 = 3.484196000000005
$ws.Range("O14").Value2 = 1.806054000000002
$ws.Range("P14").Value2 = 1.352799500000001

# --- Row 15: Rotation-NoTilt (shifted from old row 11) ---------------------
$ws.Range("B15").Value2 = "Rotation-NoTilt"
$ws.Range("C15").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 7.650550000000003
$ws.Range("F15").Value2 = 0.01
$ws.Range("G15").Value2 = 0
$ws.Range("H15").Value2 = 0
$ws.Range("I15").Value2 = 2.397512500000005
$ws.Range("J15").Value2 = 0.8908250000000003
$ws.Range("K15").Value2 = 0.04
$ws.Range("L15").Value2 = 0
$ws.Range("M15").Value2 = 0
$ws.Range("N15").Value2 = 3.825275000000001
$ws.Range("O15").Value2 = 1.915137500000001
$ws.Range("P15").Value2 = 1.373610937500001

# --- Row 16: Rotation-60detTilt (shifted from old row 12) ------------------
$ws.Range("B16").Value2 = "Rotation-60detTilt"
$ws.Range("C16").Value2 = 0.424147589324806
$ws.Range("D16").Value2 = 0.4300469861376049
$ws.Range("E16").Value2 = 4.725924619264004
$ws.Range("F16").Value2 = 0.407864331468798
$ws.Range("G16").Value2 = 0.424147589324806
$ws.Range("H16").Value2 = 0.4300469861376049
$ws.Range("I16").Value2 = 1.811344470732795
$ws.Range("J16").Value2 = 0.9468965909503991
$ws.Range("K16").Value2 = 0.4404363624447993
$ws.Range("L16").Value2 = 0.4119654467584001
$ws.Range("M16").Value2 = 0.424173062553606
$ws.Range("N16").Value2 = 2.577985802700804
$ws.Range("O16").Value2 = 1.496995881548803
$ws.Range("P16").Value2 = 1.199828299635201

# --- New row 17: HexGrid-90degTilt5degRes (shifted from old row 13) --------
# Prime formatting for the new rows by copying row 16's formats down first.
$ws.Range("A16:P16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value2 = 0.9813296137166525
$ws.Range("D17").Value2 = 0.9838600376175911
$ws.Range("E17").Value2 = 1.005877367761607
$ws.Range("F17").Value2 = 0.996241027770477
$ws.Range("G17").Value2 = 0.9813296137166525
$ws.Range("H17").Value2 = 0.9838600376175911
$ws.Range("I17").Value2 = 0.9969457478229776
$ws.Range("J17").Value2 = 0.9923119986392369
$ws.Range("K17").Value2 = 0.9959362005835859
$ws.Range("L17").Value2 = 0.9903623239135438
$ws.Range("M17").Value2 = 0.9813328007554135
$ws.Range("N17").Value2 = 0.9948687026895993
$ws.Range("O17").Value2 = 0.9918270117165819
$ws.Range("P17").Value2 = 0.992858039728209

# --- New row 18: HexGrid-90degTilt22p5degRes (shifted from old row 14) -----
$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value2 = 1.292504258358423
$ws.Range("D18").Value2 = 1.377496938145186
$ws.Range("E18").Value2 = 0.9562114723102951
$ws.Range("F18").Value2 = 0.964342631762048
$ws.Range("G18").Value2 = 1.292504258358423
$ws.Range("H18").Value2 = 1.377496938145186
$ws.Range("I18").Value2 = 0.9269566650268798
$ws.Range("J18").Value2 = 1.026646801348446
$ws.Range("K18").Value2 = 0.8806293173617331
$ws.Range("L18").Value2 = 0.9519580305852714
$ws.Range("M18").Value2 = 1.292437783746311
$ws.Range("N18").Value2 = 1.166854205227741
$ws.Range("O18").Value2 = 1.147638825143988
$ws.Range("P18").Value2 = 1.047093264362286

# --- New row 19: HexGrid-60degTilt5degRes (shifted from old row 15) --------
$ws.Range("A19").Value2 = 17
$ws.Range("B19").Value2 = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value2 = 0.9980678637560858
$ws.Range("D19").Value2 = 1.239648285085894
$ws.Range("E19").Value2 = 0.9841436777953232
$ws.Range("F19").Value2 = 0.9593104618075039
$ws.Range("G19").Value2 = 0.9980678637560858
$ws.Range("H19").Value2 = 1.239648285085894
$ws.Range("I19").Value2 = 0.9644317021841152
$ws.Range("J19").Value2 = 0.9700759480451602
$ws.Range("K19").Value2 = 0.9527240021869781
$ws.Range("L19").Value2 = 1.068919226894605
$ws.Range("M19").Value2 = 0.9981394480227012
$ws.Range("N19").Value2 = 1.111895981440608
$ws.Range("O19").Value2 = 1.045292572111202
$ws.Range("P19").Value2 = 1.017165145969458
